$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q4" sheet right after "总计", by duplicating the
#        existing "2022-Q3" sheet (same column layout/styles) and overwriting
#        its data with the new quarter's numbers. ---
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($null, $totalSheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# D:G hold numeric-looking text ("2.55", not the number 2.55) in the source
# data, matching the sibling quarter sheets. Force text via a "@" format so
# the runtime doesn't silently coerce the string to a number, then restore
# the cells to the (unstyled) default so no stray style gets attached.
foreach ($addr in @("D2", "E2", "F2", "G2")) {
    $q4Sheet.Range($addr).NumberFormat = "@"
}
$q4Sheet.Range("D2").Value = "2.55"
$q4Sheet.Range("E2").Value = "97.28"
$q4Sheet.Range("F2").Value = "0.48"
$q4Sheet.Range("G2").Value = "0.0122"
foreach ($addr in @("D2", "E2", "F2", "G2")) {
    $q4Sheet.Range($addr).Style = "Normal"
}
$q4Sheet.Range("H2").Value = 2

# --- 2. Update the "总计" summary sheet: insert a new row for 2022-Q4 right
#        after the header row, pushing the existing rows down by one. ---
# Row 5 is brand new, so first clone row 4's formatting (column A carries a
# bordered/bold style) onto it before writing the shifted-down values.
$totalSheet.Range("A4").Copy($totalSheet.Range("A5"))

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 3
$totalSheet.Range("D5").Value = 0.04

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.01

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.01

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.01
